$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19 (shifts existing rows 19-139 down to 20-140)
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new daily record
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44749
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112043
$ws.Range("G19").Value = "Pepino ensalada"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19556
$ws.Range("N19").Value = "$/caja 60 unidades"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 326
$ws.Range("Q19").Value = 60
$ws.Range("R19").Value = "Hortaliza"
